# Updated symbol list on Sun Jan 22 05:54:16 UTC 2023 with GitHub Actions
# Applies refreshed price/volume figures to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.44"
$ws.Range("E2").Value = "'-0.60%"
$ws.Range("D3").Value = "'37.72"
$ws.Range("E3").Value = "'8.32%"
$ws.Range("D4").Value = "'5.010"
$ws.Range("E4").Value = "'-2.71%"
$ws.Range("D5").Value = "'0.07864"
$ws.Range("E5").Value = "'1.46%"
$ws.Range("D6").Value = "'2.205"
$ws.Range("E6").Value = "'-6.37%"
$ws.Range("D7").Value = "'8.012"
$ws.Range("E7").Value = "'0.01%"
$ws.Range("D8").Value = "'4.003"
$ws.Range("E8").Value = "'1.42%"
$ws.Range("D9").Value = "'0.9101"
$ws.Range("E9").Value = "'-2.08%"
$ws.Range("D10").Value = "'0.1877"
$ws.Range("E10").Value = "'4.52%"
$ws.Range("D11").Value = "'0.09220"
$ws.Range("E11").Value = "'-7.23%"
$ws.Range("D12").Value = "'0.08448"
$ws.Range("E12").Value = "'-2.36%"
$ws.Range("D13").Value = "'0.03515"
$ws.Range("E13").Value = "'5.92%"
$ws.Range("D14").Value = "'0.09944"
$ws.Range("E14").Value = "'0.46%"
$ws.Range("D15").Value = "'0.001476"
$ws.Range("E15").Value = "'-1.55%"
$ws.Range("D16").Value = "'0.005633"
$ws.Range("E16").Value = "'-2.02%"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'0.34%"
$ws.Range("D18").Value = "'2.060"
$ws.Range("E18").Value = "'-3.53%"
$ws.Range("E19").Value = "'2.86%"
$ws.Range("D20").Value = "'0.1299"
$ws.Range("E20").Value = "'-2.59%"
$ws.Range("D21").Value = "'4.562"
$ws.Range("E21").Value = "'4.77%"
$ws.Range("D23").Value = "'0.04644"
$ws.Range("E23").Value = "'1.42%"
$ws.Range("E24").Value = "'0.95%"
$ws.Range("D25").Value = "'0.004448"
$ws.Range("E25").Value = "'-0.31%"
$ws.Range("E26").Value = "'-0.19%"
$ws.Range("D27").Value = "'0.0004745"
$ws.Range("E27").Value = "'39.86%"
$ws.Range("D39").Value = "'0.01749"
$ws.Range("E39").Value = "'-2.01%"
$ws.Range("D40").Value = "'0.04711"
$ws.Range("E40").Value = "'-1.66%"
$ws.Range("D41").Value = "'0.007919"
$ws.Range("E41").Value = "'1.69%"
$ws.Range("D42").Value = "'0.1389"
$ws.Range("E42").Value = "'-1.54%"
$ws.Range("D43").Value = "'0.007661"
$ws.Range("E43").Value = "'9.03%"
$ws.Range("E44").Value = "'8.90%"
$ws.Range("D45").Value = "'0.01054"
$ws.Range("E45").Value = "'14.79%"
$ws.Range("D46").Value = "'0.00005977"
$ws.Range("E46").Value = "'-2.43%"
$ws.Range("D48").Value = "'8.669"
$ws.Range("E48").Value = "'183.03%"
$ws.Range("E49").Value = "'34.80%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("D51").Value = "'0.0001998"
